$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing birthday value in D2 ---
$ws.Range("D2").Value = 33219

# --- New column J (Land) data cells must use the text number format like column I (style index 3 / numFmtId 49) ---
$ws.Range("J2:J12").NumberFormat = "@"

# Header + data written in an order that reproduces the original shared-string
# allocation order: Land, DE, AR, Ersttestung, Symptome, Uebermittelung an GSA
$ws.Range("J1").Value = "Land"

$ws.Range("J2").Value = "DE"
$ws.Range("J3").Value = "DE"
$ws.Range("J4").Value = "DE"
$ws.Range("J5").Value = "DE"
$ws.Range("J6").Value = "DE"
$ws.Range("J7").Value = "DE"
$ws.Range("J8").Value = "DE"
$ws.Range("J9").Value = "DE"
$ws.Range("J10").Value = "DE"
$ws.Range("J11").Value = "DE"

$ws.Range("J12").Value = "AR"

$ws.Range("L1").Value = "Ersttestung"
$ws.Range("M1").Value = "Symptome"
$ws.Range("K1").Value = "Übermittelung an GSA"

# --- Numeric columns K (Übermittelung an GSA), L (Ersttestung), M (Symptome) ---
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 11).Value = 1
    $ws.Cells.Item($r, 12).Value = 0
    $ws.Cells.Item($r, 13).Value = 0
}

# --- Column widths (closest achievable values given engine's width quantization) ---
$ws.Columns.Item(10).ColumnWidth = 17.333333333333336
$ws.Columns.Item(11).ColumnWidth = 18.5
$ws.Columns.Item(12).ColumnWidth = 9.333333333333332
$ws.Columns.Item(13).ColumnWidth = 8.833333333333332

# --- Selection / active cell ---
$null = $ws.Range("B11").Select()
